$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = -0.0003960461763199419
$ws.Range("G2").Value = -0.0000000391155481338501
$ws.Range("J2").Value = -0.0000000391155481338501
$ws.Range("K2").Value = 0.0003959865716751665
$ws.Range("L2").Value = 0.1771626129571466
$ws.Range("E3").Value = -0.0003960461763199419
$ws.Range("G3").Value = 0.0000000391155481338501
$ws.Range("J3").Value = 0.0000000391155481338501
$ws.Range("K3").Value = 0.0003961057809647173
$ws.Range("L3").Value = 0.1771626132236384
$ws.Range("E4").Value = 118.8882751464844
$ws.Range("F4").Value = 317.3644409179688
$ws.Range("G4").Value = 390.3087463378906
$ws.Range("H4").Value = 115.6178131103516
$ws.Range("I4").Value = -4.13555908203125
$ws.Range("J4").Value = 4.459381103515625
$ws.Range("K4").Value = -3.270462036132812
$ws.Range("L4").Value = 6.905421831948384
$ws.Range("E5").Value = 118.8882751464844
$ws.Range("F5").Value = -327.6585388183594
$ws.Range("G5").Value = 388.4421691894531
$ws.Range("H5").Value = 115.3962173461914
$ws.Range("I5").Value = -6.158538818359375
$ws.Range("J5").Value = 2.592803955078125
$ws.Range("K5").Value = -3.492057800292969
$ws.Range("L5").Value = 7.539542453444707
$ws.Range("E6").Value = 118.8882751464844
$ws.Range("F6").Value = 327.936279296875
$ws.Range("G6").Value = -387.8245239257812
$ws.Range("H6").Value = 115.0099639892578
$ws.Range("I6").Value = 6.436279296875
$ws.Range("J6").Value = -1.97515869140625
$ws.Range("K6").Value = -3.878311157226562
$ws.Range("L6").Value = 7.769700153537936
$ws.Range("E7").Value = 118.8882751464844
$ws.Range("G7").Value = -392.6463012695312
$ws.Range("H7").Value = 116.0939254760742
$ws.Range("J7").Value = -6.79693603515625
$ws.Range("K7").Value = -2.794349670410156
$ws.Range("L7").Value = 7.993446702542722
$ws.Range("E8").Value = 125.0787048339844
$ws.Range("F8").Value = 177.8002014160156
$ws.Range("G8").Value = 2.23306131362915
$ws.Range("H8").Value = 120.6702423095703
$ws.Range("I8").Value = 4.531631469726562
$ws.Range("J8").Value = 2.23306131362915
$ws.Range("K8").Value = -4.408462524414062
$ws.Range("L8").Value = 6.704982359186822
$ws.Range("E9").Value = 125.0787048339844
$ws.Range("F9").Value = -176.8183288574219
$ws.Range("G9").Value = -1.649347305297852
$ws.Range("H9").Value = 122.2574615478516
$ws.Range("I9").Value = -3.549758911132812
$ws.Range("J9").Value = -1.649347305297852
$ws.Range("K9").Value = -2.821243286132812
$ws.Range("L9").Value = 4.824992076699149
$ws.Range("E10").Value = -6.126372814178467
$ws.Range("F10").Value = -430.3816033986932
$ws.Range("G10").Value = 382.1223640173984
$ws.Range("H10").Value = -8.133339655480029
$ws.Range("I10").Value = -7.113034304841221
$ws.Range("J10").Value = 0.285755130679604
$ws.Range("K10").Value = -2.006966841301562
$ws.Range("L10").Value = 7.396271284819245
$ws.Range("E11").Value = -6.126372814178467
$ws.Range("F11").Value = -502.9451497231984
$ws.Range("G11").Value = 382.1223640173984
$ws.Range("H11").Value = -8.133339655480029
$ws.Range("I11").Value = -8.176580629346404
$ws.Range("J11").Value = 0.285755130679604
$ws.Range("K11").Value = -2.006966841301562
$ws.Range("L11").Value = 8.42413453625929
$ws.Range("E12").Value = 125.0787048339844
$ws.Range("F12").Value = -502.9451497231984
$ws.Range("G12").Value = -3.364271165890906
$ws.Range("H12").Value = 120.127938124037
$ws.Range("I12").Value = -8.176595461916918
$ws.Range("J12").Value = -3.364263536496375
$ws.Range("K12").Value = -4.950766709947402
$ws.Range("L12").Value = 10.13336437256467
$ws.Range("E13").Value = -6.126372814178467
$ws.Range("F13").Value = -423.4817509996568
$ws.Range("G13").Value = -390.7817973166199
$ws.Range("H13").Value = -5.452928036183479
$ws.Range("I13").Value = -0.2131819058047881
$ws.Range("J13").Value = -8.945188429901179
$ws.Range("K13").Value = 0.6734447779949875
$ws.Range("L13").Value = 8.973035742735522
$ws.Range("E14").Value = -6.126372814178467
$ws.Range("F14").Value = -502.9451497231984
$ws.Range("G14").Value = -390.7817973166199
$ws.Range("H14").Value = -5.452928036183479
$ws.Range("I14").Value = -8.176580629346404
$ws.Range("J14").Value = -8.945188429901179
$ws.Range("K14").Value = 0.6734447779949875
$ws.Range("L14").Value = 12.13780848026733
$ws.Range("E15").Value = -6.126372814178467
$ws.Range("F15").Value = 503.1871221643107
$ws.Range("G15").Value = 391.431068212491
$ws.Range("H15").Value = -8.887196065851157
$ws.Range("I15").Value = 8.418553070458699
$ws.Range("J15").Value = 9.594459325772277
$ws.Range("K15").Value = -2.76082325167269
$ws.Range("L15").Value = 13.05939625637436
$ws.Range("E16").Value = -6.126372814178467
$ws.Range("F16").Value = 423.3945912686113
$ws.Range("G16").Value = 391.431068212491
$ws.Range("H16").Value = -8.887196065851157
$ws.Range("I16").Value = 0.1260221747592709
$ws.Range("J16").Value = 9.594459325772277
$ws.Range("K16").Value = -2.76082325167269
$ws.Range("L16").Value = 9.984571917183342
$ws.Range("E17").Value = 125.0787048339844
$ws.Range("F17").Value = 503.1871221643107
$ws.Range("G17").Value = 8.823901994241472
$ws.Range("H17").Value = 119.4525165100649
$ws.Range("I17").Value = 8.418567903029214
$ws.Range("J17").Value = 8.823909623636004
$ws.Range("K17").Value = -5.626188323919507
$ws.Range("L17").Value = 13.43084739099577
$ws.Range("E18").Value = -6.126372814178467
$ws.Range("F18").Value = 503.1871221643107
$ws.Range("G18").Value = -379.9099285136609
$ws.Range("H18").Value = -3.151352475130295
$ws.Range("I18").Value = 8.418553070458699
$ws.Range("J18").Value = 1.926680373057877
$ws.Range("K18").Value = 2.975020339048172
$ws.Range("L18").Value = 9.134269487912338
$ws.Range("E19").Value = -6.126372814178467
$ws.Range("F19").Value = 437.0525198116447
$ws.Range("G19").Value = -379.9099285136609
$ws.Range("H19").Value = -3.151352475130295
$ws.Range("I19").Value = 13.78395071779272
$ws.Range("J19").Value = 1.926680373057877
$ws.Range("K19").Value = 2.975020339048172
$ws.Range("L19").Value = 14.2323624415701
